$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4   = -20.178
    7   = -19.851
    16  = -21.901
    28  = -21.899
    29  = -21.344
    32  = -21.816
    40  = -19.965
    52  = -22.058
    57  = -22.275
    66  = -21.616
    100 = -22.029
}

foreach ($row in $updates.Keys) {
    $ws.Range("A$row").Value = $updates[$row]
}
